# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns,
# plus two rows (39/40) whose Coin/Link/Price/Volume all changed.
# Note: price strings that look like plain decimals (e.g. "242.26") are
# entered with a leading apostrophe so Excel keeps them as text, matching
# the existing text-typed Price column (values with two dots, like
# "29.391.40", already stay text on their own).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.391.40"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.875.64"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'0.7127"
$ws.Range("D6").Value = "'242.26"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.3117"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "'0.07782"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'25.12"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").Value = "1.866.76"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "'5.244"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'0.7131"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "'91.26"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "29.386.64"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'0.000008271"
$ws.Range("E17").Value = "  +5.99%  "
$ws.Range("D18").Value = "'6.050"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'241.47"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "'13.27"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "2.121.91"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D23").Value = "'7.797"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'0.1606"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "'164.10"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").Value = "'9.074"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").Value = "'18.51"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").Value = "'1.511"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "'4.433"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'4.322"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").Value = "'1.288"
$ws.Range("E32").Value = "  -4.04%  "
$ws.Range("D33").Value = "'0.05309"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("D34").Value = "'1.943"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "'1.180"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "'0.7483"
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("D37").Value = "'2.696"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "'0.01872"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.208.53"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.723"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").Value = "'6.452"
$ws.Range("E41").Value = "  +3.50%  "
$ws.Range("D42").Value = "'0.8883"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").Value = "'72.89"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'109.57"
$ws.Range("E44").Value = "  +7.28%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "2.019.14"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "'1.822"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("D48").Value = "'0.5211"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  +5.79%  "
$ws.Range("D50").Value = "'9.388"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "'0.4324"
$ws.Range("E51").Value = "  +1.10%  "
